$d = $word.ActiveDocument

# The paragraph "Nossa música hoje  >   Conteúdo da tabela" loses its trailing
# "Conteúdo da tabela" text, and the hidden "_GoBack" bookmark (currently sitting
# at the end of the "Sertanejo" paragraph, right after its final "OK") moves to
# sit at the end of that now-shorter "Nossa música hoje" paragraph.

# Locate "Conteúdo da tabela" (find only, no replace yet) and collapse the
# match range to its start - that position is exactly where the paragraph will
# end once the text is removed, i.e. where the bookmark needs to live.
$markRange = $d.Content
$markRange.Find.Execute("Conteúdo da tabela") | Out-Null
$markRange.Collapse(1)

# Remove the old "_GoBack" bookmark from the end of the "Sertanejo" paragraph.
$d.Bookmarks("_GoBack").Delete()

# Re-create "_GoBack" at the new location now, before deleting the text
# (adding it after the shrink would land on a stale/degenerate range).
$d.Bookmarks.Add("_GoBack", $markRange)

# Finally, delete the "Conteúdo da tabela" text that used to close the
# "Nossa música hoje" paragraph.
$d.Content.Find.Execute("Conteúdo da tabela", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
